{"js": "// The edit replaces the text of every run in this document (the header date\n// plus all 25 \"two-digit x two-digit\" answer cells) with a new value, one-\n// for-one, in document order. Every old value below is unique in the body,\n// so a literal (non-wildcard) search-and-replace for each pair is safe and\n// unambiguous.\nconst body = context.document.body;\n\nconst replacements = [\n  [\n    \"2024-08-27 Tuesday\",\n    \"2024-08-28 Wednesday\"\n  ],\n  [\n    \"88\u00d720=1760\",\n    \"56\u00d719=1064\"\n  ],\n  [\n    \"47\u00d763=2961\",\n    \"14\u00d721=294\"\n  ],\n  [\n    \"66\u00d760=3960\",\n    \"72\u00d738=2736\"\n  ],\n  [\n    \"19\u00d757=1083\",\n    \"36\u00d782=2952\"\n  ],\n  [\n    \"60\u00d785=5100\",\n    \"11\u00d744=484\"\n  ],\n  [\n    \"80\u00d726=2080\",\n    \"17\u00d785=1445\"\n  ],\n  [\n    \"21\u00d711=231\",\n    \"14\u00d744=616\"\n  ],\n  [\n    \"17\u00d735=595\",\n    \"91\u00d740=3640\"\n  ],\n  [\n    \"91\u00d760=5460\",\n    \"62\u00d732=1984\"\n  ],\n  [\n    \"61\u00d722=1342\",\n    \"35\u00d721=735\"\n  ],\n  [\n    \"95\u00d785=8075\",\n    \"90\u00d723=2070\"\n  ],\n  [\n    \"34\u00d766=2244\",\n    \"75\u00d735=2625\"\n  ],\n  [\n    \"23\u00d788=2024\",\n    \"85\u00d785=7225\"\n  ],\n  [\n    \"91\u00d785=7735\",\n    \"20\u00d796=1920\"\n  ],\n  [\n    \"47\u00d712=564\",\n    \"98\u00d766=6468\"\n  ],\n  [\n    \"75\u00d733=2475\",\n    \"93\u00d798=9114\"\n  ],\n  [\n    \"21\u00d741=861\",\n    \"12\u00d734=408\"\n  ],\n  [\n    \"64\u00d726=1664\",\n    \"48\u00d761=2928\"\n  ],\n  [\n    \"45\u00d771=3195\",\n    \"20\u00d793=1860\"\n  ],\n  [\n    \"99\u00d742=4158\",\n    \"47\u00d723=1081\"\n  ],\n  [\n    \"86\u00d785=7310\",\n    \"43\u00d741=1763\"\n  ],\n  [\n    \"37\u00d745=1665\",\n    \"39\u00d797=3783\"\n  ],\n  [\n    \"58\u00d719=1102\",\n    \"13\u00d780=1040\"\n  ],\n  [\n    \"92\u00d713=1196\",\n    \"62\u00d721=1302\"\n  ],\n  [\n    \"16\u00d746=736\",\n    \"28\u00d761=1708\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `expected exactly 1 match for \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The edit replaces the text of every run in this document (the header date\n# plus all 25 \"two-digit x two-digit\" answer cells) with a new value, one-for-\n# one, in document order. Every old value below is unique in the body, so a\n# literal (non-wildcard) Find/Replace for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-08-27 Tuesday\", \"2024-08-28 Wednesday\"),\n  @(\"88\u00d720=1760\", \"56\u00d719=1064\"),\n  @(\"47\u00d763=2961\", \"14\u00d721=294\"),\n  @(\"66\u00d760=3960\", \"72\u00d738=2736\"),\n  @(\"19\u00d757=1083\", \"36\u00d782=2952\"),\n  @(\"60\u00d785=5100\", \"11\u00d744=484\"),\n  @(\"80\u00d726=2080\", \"17\u00d785=1445\"),\n  @(\"21\u00d711=231\", \"14\u00d744=616\"),\n  @(\"17\u00d735=595\", \"91\u00d740=3640\"),\n  @(\"91\u00d760=5460\", \"62\u00d732=1984\"),\n  @(\"61\u00d722=1342\", \"35\u00d721=735\"),\n  @(\"95\u00d785=8075\", \"90\u00d723=2070\"),\n  @(\"34\u00d766=2244\", \"75\u00d735=2625\"),\n  @(\"23\u00d788=2024\", \"85\u00d785=7225\"),\n  @(\"91\u00d785=7735\", \"20\u00d796=1920\"),\n  @(\"47\u00d712=564\", \"98\u00d766=6468\"),\n  @(\"75\u00d733=2475\", \"93\u00d798=9114\"),\n  @(\"21\u00d741=861\", \"12\u00d734=408\"),\n  @(\"64\u00d726=1664\", \"48\u00d761=2928\"),\n  @(\"45\u00d771=3195\", \"20\u00d793=1860\"),\n  @(\"99\u00d742=4158\", \"47\u00d723=1081\"),\n  @(\"86\u00d785=7310\", \"43\u00d741=1763\"),\n  @(\"37\u00d745=1665\", \"39\u00d797=3783\"),\n  @(\"58\u00d719=1102\", \"13\u00d780=1040\"),\n  @(\"92\u00d713=1196\", \"62\u00d721=1302\"),\n  @(\"16\u00d746=736\", \"28\u00d761=1708\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Text = $oldText\n  $range.Find.Replacement.Text = $newText\n  $range.Find.MatchCase = $true\n  $range.Find.MatchWildcards = $false\n\n  # wdReplaceOne (1) with Forward=$true / Wrap=wdFindContinue (1): each old value\n  # occurs exactly once, so replace that single occurrence.\n  $found = $range.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 1)\n  if (-not $found) {\n    throw \"Find/Replace did not find expected text: $oldText\"\n  }\n}\n"}
